# Scientific Writing Checklist - apply commit:
# "Added comments on what should be in the background chapter + on good
#  scientific figures + on how to write clearly and concisely."
#
# Concretely (per the OOXML diff) this edit:
#  1. Splits the single "Are all figures formatted nicely..." bullet into a
#     short top-level bullet plus four new second-level sub-bullets.
#  2. Expands the "Are all text elements in figures readable?" bullet with
#     an extra (underlined) clause about consistent figure text size.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Split "Are all figures formatted nicely and visually appealing? (No
#    bad scans, badly scaled screenshots, etc. - generally, prefer vector
#    graphics wherever possible.)" into a short bullet + 4 sub-bullets.
# ---------------------------------------------------------------------

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Are all figures formatted nicely")) {
        $target = $i
        break
    }
}

$p = $d.Paragraphs($target)
$rng = $d.Range($p.Range.Start, $p.Range.End - 1)
$rng.Text = "Are all figures formatted nicely and visually appealing?"

# Insert four empty paragraphs right after it (they inherit the List
# Paragraph style / numId=3 / ilvl=0 numbering of their predecessor).
$anchor = $d.Paragraphs($target)
$anchor.Range.InsertParagraphAfter()
$anchor.Range.InsertParagraphAfter()
$anchor.Range.InsertParagraphAfter()
$anchor.Range.InsertParagraphAfter()

$subTexts = @(
    "No bad scans, badly scaled screenshots, etc. " + [char]0x2013 + " generally, prefer vector graphics wherever possible.",
    "Not too large (wasting a lot of space, looking disproportionate), not too small (hard to recognize details). This also applies to line widths, which should neither be too thin nor too thick.",
    "Visually appealing and clearly distinguishable color map used throughout the thesis.",
    "No unnecessary details that are not relevant for the information you wish to convey."
)

for ($k = 0; $k -lt 4; $k++) {
    $sp = $d.Paragraphs($target + 1 + $k)
    $sp.Range.Text = $subTexts[$k]
    # Demote to second list level (ilvl=1) and apply the explicit
    # indentation Word uses for that level.
    $sp.Range.ListFormat.ListLevelNumber = 2
    $sp.LeftIndent = 49.5
    $sp.FirstLineIndent = -13.5
}

# Italicize the word "and" in the color-map bullet ("Visually appealing
# *and* clearly distinguishable color map ...").
$colorMapPara = $d.Paragraphs($target + 3)
$findRange = $colorMapPara.Range
$findRange.Find.Execute("and", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$findRange.Italic = 1

# ---------------------------------------------------------------------
# 2. Expand the figure-text-size bullet with the extra underlined clause.
# ---------------------------------------------------------------------

$target2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Are all text elements in figures readable?")) {
        $target2 = $i
        break
    }
}

$p2 = $d.Paragraphs($target2)
$rng2 = $d.Range($p2.Range.Start, $p2.Range.End - 1)
$rng2.Text = "Are all text elements in figures readable? (Text size in figures generally should be consistent across figures and should be about 1pt smaller than normal text size in the document.)"

$underlineRange = $d.Paragraphs($target2).Range
$underlineRange.Find.Execute("Text size in figures generally should be consistent across figures and should be about 1pt smaller than normal text size", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$underlineRange.Underline = 1

Write-Output "Done."
